# Improved accuracy of cosh function.
# Update the underlying measurement data on Sheet1 (test/results.xlsx) for
# rows 24 (sinh), 25 (cosh), 26 (tanh) and 29 (atanh). The N/O/V/W columns
# are computed by shared formulas (L/F, M/G, T/F, U/G respectively) and the
# summary rows 33/34 are AVERAGE() formulas, so they recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24 - sinh
$ws.Range("F24").Value = 8956
$ws.Range("G24").Value = 8358
$ws.Range("L24").Value = 10288
$ws.Range("M24").Value = 3095
$ws.Range("T24").Value = 10423
$ws.Range("U24").Value = 2961

# Row 25 - cosh
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 7230
$ws.Range("G25").Value = 6475
$ws.Range("L25").Value = 10276
$ws.Range("M25").Value = 3171
$ws.Range("T25").Value = 10376
$ws.Range("U25").Value = 2925

# Row 26 - tanh
$ws.Range("F26").Value = 9782
$ws.Range("G26").Value = 8505
$ws.Range("L26").Value = 4959
$ws.Range("M26").Value = 1667
$ws.Range("T26").Value = 5029
$ws.Range("U26").Value = 1626

# Row 29 - atanh
$ws.Range("F29").Value = 3975
$ws.Range("G29").Value = 3351
$ws.Range("L29").Value = 1893
$ws.Range("M29").Value = 1233
$ws.Range("T29").Value = 1916
$ws.Range("U29").Value = 1358

$excel.CalculateFullRebuild()
